$d = $word.ActiveDocument

# Step 1: remove the _GoBack bookmark in paragraph 1
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}

# Step 2: delete old paragraph 6 (space) and paragraph 5 (old Figure 3 "plan view" caption)
$d.Paragraphs(6).Range.Delete()
$d.Paragraphs(5).Range.Delete()

# Now paragraphs are:
# 1: Figure 1 caption
# 2: (space) -> becomes Figure 2 caption
# 3: old Figure 2 text -> becomes Figure 3 caption
# 4: (space) -> becomes Figures 4a and b caption
# 5: old Figure 4a and b caption -> to be deleted (content merged into 4)

# Step 3: set paragraph 2 text (new Figure 2 caption)
$r2 = $d.Paragraphs(2).Range
$r2.MoveEnd(1, -1)
$r2.Text = "Figure 2: A simplified beach ecosystem food web graphic showing select species, their relative trophic positions, and common habitats along the beach. The graph shown here is specific to the various species mentioned in the paper."

# Step 4: set paragraph 3 text (new Figure 3 caption)
$r3 = $d.Paragraphs(3).Range
$r3.MoveEnd(1, -1)
$r3.Text = "Figure 3: Graphic showing beach in elevation at three stages in nourishment evolution. Profile A shows an idealized typical beach cross-shore profile prior to nourishment. Profile B depicts the beach profile, with the new sand volume added in the immediate aftermath of nourishment. The placed fill sand volume is positioned principally on the subaerial beach and is bulldozed into an initial design profile from which natural processes driven by local winds and waves will continue with redistribution toward equilibrium with local conditions. Profile C shows the profile after placement at point where sand has been naturally redistributed onshore and offshore toward a profile geometry that is equilibrated to local conditions. "

# Step 5: set paragraph 4 text (new Figures 4a and b caption)
$r4 = $d.Paragraphs(4).Range
$r4.MoveEnd(1, -1)
$r4.Text = "Figures 4a and b: Typical beach scarps formed when the beach profile is out of equilibrium with current wind and wave conditions--in such cases wave energies in the adjacent surf zone are sufficiently high to erode sands from the subaerial beach, moving them offshore leaving behind steeply sloped scarp features along the beach. Scarp elevations can range from a few centimeters to three meters or more. Such features are commonplace under erosive conditions such as during storms. They are also common following nourishment sand placement as the unstable beach profile undergoes morphological adjustment.  "

# Step 6: delete old paragraph 5 (duplicate Figure 4a and b text, now redundant)
$d.Paragraphs(5).Range.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host "=== Para $i ==="
    Write-Host $d.Paragraphs($i).Range.Text
}
